$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("rec1", "cnumani@outlook.com", "2024-10-17"),
    @("rec2", "reddy.binary@gmail.com", "2024-10-17"),
    @("rec3", "indumathilalam03@gmail.com", "2024-10-17"),
    @("rec1", "cnumani@outlook.com", "2024-10-17"),
    @("rec2", "reddy.binary@gmail.com", "2024-10-17"),
    @("rec3", "indumathilalam03@gmail.com", "2024-10-17"),
    @("rec1", "cnumani@outlook.com", "2024-10-17"),
    @("rec2", "reddy.binary@gmail.com", "2024-10-17"),
    @("rec3", "indumathilalam03@gmail.com", "2024-10-17"),
    @("cnumani", "cnumani@outlook.com", "2024-10-17"),
    @("binary", "reddy.binary@gmail.com", "2024-10-17"),
    @("indumanthi", "indumathilalam03@gmail.com", "2024-10-17")
)

$row = 4
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $dateCell = $ws.Cells.Item($row, 3)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $entry[2]
    $dateCell.ClearFormats()
    $row++
}
